$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so values like "1.000" or "240.67" keep
# their literal formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 2-20: only the Price (D) and Volume(1h) (E) columns change this refresh
$ws.Range("D2").Value = "29.903.11"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.902.94"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "0.8023"
$ws.Range("E5").Value = "  +5.89%  "
$ws.Range("D6").Value = "240.67"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.3107"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "26.37"
$ws.Range("E9").Value = "  +3.70%  "
$ws.Range("D10").Value = "0.06980"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("D11").Value = "0.07987"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.917.54"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "0.7376"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "5.157"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "92.22"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "29.919.71"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "13.92"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "5.841"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "244.01"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "0.000007769"
$ws.Range("E20").Value = "  +0.73%  "

# A new coin (WrappedliquidstakedEther2.0) enters the ranked list at row 21,
# shifting every following coin down by one row; the former last row (Cronos)
# drops off the bottom of the table entirely.
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.157.09"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.881"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "167.48"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.166"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.1451"
$ws.Range("E27").Value = "  +14.20%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "18.83"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.056"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.353"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.510"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.273"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05504"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.049"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.259"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7296"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01921"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.787"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.4390"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "72.08"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.959"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8356"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "1.876"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "100.56"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.519"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.648"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.067.86"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "976.09"
$ws.Range("E50").Value = "  +7.73%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "36.07"
$ws.Range("E51").Value = "  -0.46%  "
